$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: add "client" / "clientemail" headers (keeps existing cell styles s="4"/s="5")
$ws.Range("A2").Value = "client"
$ws.Range("B2").Value = "clientemail"

# Row 10: add an empty, styled C10 cell (matches D10's font/style, s="3") next to D10
$ws.Range("C10").Font.Underline = $true

# Move the active selection to B2
$ws.Range("B2").Select() | Out-Null
